$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("weeknr 45")

# Row 8: played the game PyramidPanic from 09:35 to 10:00
$ws.Range("C8").Value = 0.39930555555555558
$ws.Range("D8").Value = 0.41666666666666669
$ws.Range("F8").Value = "Het spel PyramidPanic gespeeld."

# Row 9: icon added, name/assets work, from 10:00 to 10:17
$ws.Range("C9").Value = 0.41666666666666669
$ws.Range("D9").Value = 0.4284722222222222
$ws.Range("F9").Value = "Ico toe gevoegd. Naam in PyramidPanic gemaakt. En assets toegevoegd aan game"
$ws.Rows.Item(9).RowHeight = 30

# Update the selected cell shown in the saved view
$ws.Activate()
$ws.Range("C10").Select()
